$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Add new header "Servised by" in O1, matching the header style used by the
# rest of row 1 (bold font + border + centered alignment == existing style
# of N1), without introducing a brand-new style entry.
$ws.Cells.Item(1, 15).Value = "Servised by"
$ws.Cells.Item(1, 14).Copy()
$ws.Cells.Item(1, 15).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Extend the new column down through the data rows (2-13) as blank cells so
# the column exists across the whole table / used range.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 15).Style = "Normal"
}
